# Apply updated cryptocurrency price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing text storage so
# numeric-looking strings (e.g. "305.50") are not auto-coerced into
# numbers (which would drop the formatting / trailing zeros).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue $ws.Range("D2") "41.863.27"
$ws.Range("E2").Value = "  -0.43%  "

Set-TextValue $ws.Range("D3") "2.270.02"
$ws.Range("E3").Value = "  +0.31%  "

$ws.Range("E4").Value = "  -0.01%  "

Set-TextValue $ws.Range("D5") "305.50"
$ws.Range("E5").Value = "  +1.01%  "

Set-TextValue $ws.Range("D6") "93.03"
$ws.Range("E6").Value = "  +0.26%  "

Set-TextValue $ws.Range("D7") "0.530"
$ws.Range("E7").Value = "  -0.42%  "

$ws.Range("E8").Value = "  -0.05%  "

Set-TextValue $ws.Range("D9") "0.486"
$ws.Range("E9").Value = "  +0.39%  "

Set-TextValue $ws.Range("D10") "32.70"
$ws.Range("E10").Value = "  -0.13%  "

Set-TextValue $ws.Range("D11") "0.0798"
$ws.Range("E11").Value = "  -0.12%  "

$ws.Range("E12").Value = "  -1.92%  "

Set-TextValue $ws.Range("D13") "6.68"
$ws.Range("E13").Value = "  -0.12%  "

Set-TextValue $ws.Range("D14") "2.619.26"
$ws.Range("E14").Value = "  +0.27%  "

Set-TextValue $ws.Range("D15") "14.34"
$ws.Range("E15").Value = "  +1.44%  "

Set-TextValue $ws.Range("D16") "2.268.09"
$ws.Range("E16").Value = "  +0.12%  "

Set-TextValue $ws.Range("D17") "0.783"
$ws.Range("E17").Value = "  +3.29%  "

Set-TextValue $ws.Range("D18") "41.782.68"
$ws.Range("E18").Value = "  -0.35%  "

$ws.Range("E19").Value = "  +6.07%  "

Set-TextValue $ws.Range("D20") "0.0₃0920"
$ws.Range("E20").Value = "  +1.37%  "

Set-TextValue $ws.Range("D21") "5.98"
$ws.Range("E21").Value = "  +0.77%  "

Set-TextValue $ws.Range("D22") "68.07"
$ws.Range("E22").Value = "  +1.13%  "

Set-TextValue $ws.Range("D23") "244.24"
$ws.Range("E23").Value = "  +1.03%  "

Set-TextValue $ws.Range("D24") "2.59"
$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("E25").Value = "  +1.81%  "

$ws.Range("E26").Value = "  -0.02%  "

Set-TextValue $ws.Range("D27") "24.03"
$ws.Range("E27").Value = "  +0.41%  "

Set-TextValue $ws.Range("D28") "9.67"
$ws.Range("E28").Value = "  -0.32%  "

Set-TextValue $ws.Range("D29") "2.08"
$ws.Range("E29").Value = "  -5.03%  "

Set-TextValue $ws.Range("D30") "34.82"
$ws.Range("E30").Value = "  +1.75%  "

Set-TextValue $ws.Range("D31") "159.34"
$ws.Range("E31").Value = "  +0.38%  "

Set-TextValue $ws.Range("D32") "5.35"
$ws.Range("E32").Value = "  +3.53%  "

$ws.Range("E33").Value = "  +0.00%  "

Set-TextValue $ws.Range("D34") "0.0744"
$ws.Range("E34").Value = "  +0.20%  "

Set-TextValue $ws.Range("D35") "3.03"
$ws.Range("E35").Value = "  -1.56%  "

Set-TextValue $ws.Range("D36") "17.09"
$ws.Range("E36").Value = "  +2.83%  "

$ws.Range("E37").Value = "  -1.47%  "

$ws.Range("E38").Value = "  +0.96%  "

$ws.Range("E39").Value = "  +0.70%  "

$ws.Range("E40").Value = "  -0.24%  "

Set-TextValue $ws.Range("D41") "3.94"
$ws.Range("E41").Value = "  -0.39%  "

Set-TextValue $ws.Range("D44") "2.26"
$ws.Range("E44").Value = "  +12.80%  "

Set-TextValue $ws.Range("D45") "0.0283"
$ws.Range("E45").Value = "  +1.19%  "

$ws.Range("E46").Value = "  +1.55%  "

$ws.Range("E47").Value = "  +0.32%  "

Set-TextValue $ws.Range("D48") "53.59"
$ws.Range("E48").Value = "  +3.28%  "

Set-TextValue $ws.Range("D49") "73.24"
$ws.Range("E49").Value = "  +3.66%  "

Set-TextValue $ws.Range("D50") "1.51"
$ws.Range("E50").Value = "  -0.83%  "

Set-TextValue $ws.Range("D51") "1.15"
$ws.Range("E51").Value = "  +0.09%  "

# Row 42/43: ranking changed places, swapping Maker <-> EnergySwap
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D42") "19.77"
$ws.Range("E42").Value = "  -1.53%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D43") "2.013.38"
$ws.Range("E43").Value = "  -2.11%  "
